$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The cheat-sheet table used to have its "Save File / Save As / Undo / Redo /
# New Project / Open Project" rows filed under the "General" category. This
# upload splits those six rows out into their own new "File" category (placed
# first) and moves the "Video" category (previously first) down to the very
# bottom of the list. Every other category keeps its internal row order.
$data = @(
    @("File", "Save File", "Ctrl + S"),
    @("File", "Save As", "Ctrl + Shift + S"),
    @("File", "Undo", "Ctrl + Z"),
    @("File", "Redo", "Ctrl + Shift + Z"),
    @("File", "New Project", "Ctrl + N"),
    @("File", "Open Project", "Ctrl + O"),
    @("General", "Delete", "Delete"),
    @("General", "Focus on Selected Object", "Z"),
    @("General", "Duplicate Linked", "Ctrl + D"),
    @("General", "Duplicate Unique", "Alt + D"),
    @("General", "Select All", "Alt + A"),
    @("General", "Group", "Ctrl + G"),
    @("General", "Ungroup", "Ctrl + Shift + G"),
    @("General", "Lock / Unlock Object", "Ctrl + L"),
    @("General", "Show / Hide", "Ctrl + H"),
    @("View", "Switch Orbit / Fly Mode", "Alt + V"),
    @("View", "Enter / Exit Walk Mode", "Alt + W"),
    @("View", "Rotate View", "MMB / RMB"),
    @("View", "Pan View", "Shift + MMB / RMB"),
    @("View", "Move Forward", "W"),
    @("View", "Move Backward", "S"),
    @("View", "Move Left", "A"),
    @("View", "Move Right", "D"),
    @("View", "Move Up", "Q"),
    @("View", "Move Down", "E"),
    @("Insert", "Insert Point Light", "'1"),
    @("Insert", "Insert Spot Light", "'2"),
    @("Insert", "Insert Strip Light", "'3"),
    @("Insert", "Insert Rectangle Light", "'4"),
    @("Tool", "Switch Move / Rotate / Scale", "V"),
    @("Tool", "Material Picker", "I"),
    @("Tool", "Material Brush", "O"),
    @("Tool", "Rotate Before Placement", "R"),
    @("Tool", "Scale Before Placement", "C"),
    @("Tool", "Drop Object Vertically", "Shift + F"),
    @("Tool", "Drop Object Along Terrain", "Ctrl + F"),
    @("Tool", "Brush Tool", "B"),
    @("Tool", "Eraser Tool", "Alt"),
    @("Video", "Play / Pause", "Space"),
    @("Video", "Timeline Zoom", "Alt + Scroll"),
    @("Video", "Add Keyframe", "K"),
    @("Video", "Add Current View", "Alt + C"),
    @("Video", "Add Selected Camera View", "Shift + Alt + C")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
}

# Update the saved view state: scrolled down so row 34 is at the top, with
# B7 selected (was O3 / O1:O3 before).
$ws.Application.Goto($ws.Range("A34"), $false)
$ws.Range("B7").Select()
